$wb = $excel.ActiveWorkbook

# --- Update text on "Hoja1" (A1) with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.12 = 28490.17 pesos`n✅ 28490.17 pesos = 7.09 = 932.8 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update numeric rate cells on "tasas" sheet ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 140.4
$ws2.Range("O10").Value = 4000.02
$ws2.Range("N12").Value = 4021
$ws2.Range("O12").Value = 131.652
